$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (inline strings) and column C (numbers) for rows 2-16
$ws.Range("B2").Value = "<then>"
$ws.Range("C2").Value = 51

$ws.Range("B3").Value = "<him>"
$ws.Range("C3").Value = 55

$ws.Range("B4").Value = "<then>"
$ws.Range("C4").Value = 59

$ws.Range("B5").Value = "<has>"
$ws.Range("C5").Value = 57

$ws.Range("B6").Value = "<so>"
$ws.Range("C6").Value = 52

$ws.Range("B7").Value = "<on>"
$ws.Range("C7").Value = 51

$ws.Range("B8").Value = "<fonwa>"
$ws.Range("C8").Value = 52

$ws.Range("B9").Value = "<wound>"
$ws.Range("C9").Value = 57

$ws.Range("B10").Value = "<the>"
$ws.Range("C10").Value = 55

$ws.Range("B11").Value = "<wre>"
$ws.Range("C11").Value = 53

$ws.Range("B12").Value = "<the>"
$ws.Range("C12").Value = 56

$ws.Range("B13").Value = "<use>"

$ws.Range("B14").Value = "<they>"
$ws.Range("C14").Value = 56

$ws.Range("B15").Value = "<tine>"
$ws.Range("C15").Value = 57

$ws.Range("C16").Value = 28
